$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B13 was stored as text "3"; convert it to a real numeric value 3
$ws.Range("B13").Value = 3

# Append new row 14 with annotation data
$ws.Range("A14").Value = "Ying Tang"

# B14 must stay a text value "1" (not get auto-converted to a number)
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "1"
$ws.Range("B14").Style = "Normal"

$ws.Range("C14").Value = "absolutely know,  hide"
$ws.Range("D14").Value = "CRT"
$ws.Range("E14").Value = "RES"
$ws.Range("F14").Value = "77474e59-42ef-43e4-850b-a07d6b41a266"
$ws.Range("G14").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H14").Value = "You absolutely know this but you hide these results."
